# EXP.xlsx edits
#  1. Tweak a couple of existing RMSE-like values (L16 / L17).
#  2. Record a new "bagging" experiment row (row 23) plus its paired
#     C6H6/CO comparison block (row 24) with a tansig*2 DNN run.
#  3. Move the active selection to I24 to match where the work left off.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. corrected RMSE values for the two "S1,2,3,4" lassoCV(cv=3,normalize) rows ---
$ws.Range("L16").Value = 0.42303
$ws.Range("L17").Value = 0.42303

# --- 2. new row 23: "bagging" result block (columns I:M) + CO label (column T) ---
# Write T23 ("CO") first: it duplicates the text already used elsewhere on the
# sheet, then T24's new tansig label, then I23's new "bagging" label -- in that
# order -- so any newly-minted shared strings land in the same relative order
# the workbook ended up with.
$ws.Range("T23").Value = "CO"

$ws.Range("I23").Value = "bagging"
$ws.Range("J23").Value = 0.284
$ws.Range("K23").Value = -0.0009
$ws.Range("L23").Value = 0.4206
$ws.Range("M23").Value = 0.2298
$ws.Range("M23").NumberFormat = "0.0000_ "

# --- 3. new row 24: tansig*2 DNN run summary (columns T:X) ---
$ws.Range("T24").Value = "tansig*2,  输出不经过norm, hw=5, 8000(已收敛)"
$ws.Range("U24").Value = 0.317079656256
$ws.Range("V24").Value = 1.42910258876
$ws.Range("W24").Value = 0.115775787968
$ws.Range("X24").Value = 0.13396665
$ws.Range("X24").NumberFormat = "0.0000_ "

# --- 4. move the selection to where editing left off ---
$ws.Range("I24").Select()
